# Implemented more 'Get Bucket' VIs
# Adds a new error-code row (412348 / NoSuchWebsiteConfiguration) to the
# S3 error-codes table on Sheet1, just below the existing last row (49).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New row 50: error code number in column A, description in column B.
$ws.Range("A50").Value = 412348
$ws.Range("B50").Value = "NoSuchWebsiteConfiguration - The specified bucket does not have a website configuration."

# Move / record the active selection the same way the authored workbook
# ends up (active cell two rows below the new last data row).
$null = $ws.Range("B54").Select()
